$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-02-11 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-12 Monday", 2)

# Update the division problems in the table, addressed by row/column
# to avoid ambiguity from values that coincide with other cells' old/new text.
$tbl = $d.Tables.Item(1)
$tbl.Cell(1, 1).Range.Text = "85÷3="  # was 56÷9=
$tbl.Cell(1, 2).Range.Text = "25÷3="  # was 20÷4=
$tbl.Cell(1, 3).Range.Text = "97÷2="  # was 81÷8=
$tbl.Cell(1, 4).Range.Text = "15÷7="  # was 46÷5=
$tbl.Cell(1, 5).Range.Text = "16÷9="  # was 47÷5=
$tbl.Cell(5, 1).Range.Text = "77÷4="  # was 18÷2=
$tbl.Cell(5, 2).Range.Text = "20÷3="  # was 35÷2=
$tbl.Cell(5, 3).Range.Text = "76÷4="  # was 32÷2=
$tbl.Cell(5, 4).Range.Text = "76÷2="  # was 31÷8=
$tbl.Cell(5, 5).Range.Text = "96÷3="  # was 85÷5=
$tbl.Cell(9, 1).Range.Text = "44÷2="  # was 95÷6=
$tbl.Cell(9, 2).Range.Text = "35÷2="  # was 48÷6=
$tbl.Cell(9, 3).Range.Text = "82÷7="  # was 63÷9=
$tbl.Cell(9, 4).Range.Text = "27÷5="  # was 97÷2=
$tbl.Cell(9, 5).Range.Text = "34÷5="  # was 52÷3=
$tbl.Cell(13, 1).Range.Text = "59÷6="  # was 13÷2=
$tbl.Cell(13, 2).Range.Text = "23÷9="  # was 95÷7=
$tbl.Cell(13, 3).Range.Text = "16÷3="  # was 48÷7=
$tbl.Cell(13, 4).Range.Text = "34÷7="  # was 86÷6=
$tbl.Cell(13, 5).Range.Text = "82÷5="  # was 60÷6=
$tbl.Cell(17, 1).Range.Text = "79÷3="  # was 40÷8=
$tbl.Cell(17, 2).Range.Text = "16÷5="  # was 16÷9=
$tbl.Cell(17, 3).Range.Text = "57÷4="  # was 18÷4=
$tbl.Cell(17, 4).Range.Text = "65÷5="  # was 57÷3=
$tbl.Cell(17, 5).Range.Text = "56÷5="  # was 22÷8=
